$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '2026-02-15 17:48:34'
$ws.Range('E3').Value = '2026-02-15 17:48:36'
$ws.Range('G3').Value = '192 cm'
$ws.Range('O3').Value = '-5.5 °C'
$ws.Range('E4').Value = '2026-02-15 17:48:38'
$ws.Range('O4').Value = '7.2 °C'
$ws.Range('E5').Value = '2026-02-15 17:48:41'
$ws.Range('I5').Value = '3.2 mm'
$ws.Range('K5').Value = '5.8 MJ/m2'
$ws.Range('E6').Value = '2026-02-15 17:48:43'
$ws.Range('O6').Value = '8.4 °C'
$ws.Range('E7').Value = '2026-02-15 17:48:46'
$ws.Range('E8').Value = '2026-02-15 17:48:48'
$ws.Range('O8').Value = '7.9 °C'
$ws.Range('E9').Value = '2026-02-15 17:48:51'
$ws.Range('H9').Value = '47%'
$ws.Range('E10').Value = '2026-02-15 17:48:53'
$ws.Range('H10').Value = '68%'
$ws.Range('E11').Value = '2026-02-15 17:48:56'
$ws.Range('E12').Value = '2026-02-15 17:48:58'
$ws.Range('H12').Value = '51%'
$ws.Range('O12').Value = '10.9 °C'
$ws.Range('E13').Value = '2026-02-15 17:48:59'
$ws.Range('O13').Value = '6.4 °C'
$ws.Range('E14').Value = '2026-02-15 17:49:00'
$ws.Range('H14').Value = '56%'
$ws.Range('K14').Value = '11.6 MJ/m2'
$ws.Range('E15').Value = '2026-02-15 17:49:02'
$ws.Range('H15').Value = '47%'
$ws.Range('E16').Value = '2026-02-15 17:49:03'
$ws.Range('H16').Value = '59%'
$ws.Range('O16').Value = '-2.2 °C'
$ws.Range('E17').Value = '2026-02-15 17:49:04'
$ws.Range('H17').Value = '34%'
$ws.Range('E18').Value = '2026-02-15 17:49:05'
$ws.Range('O18').Value = '7.1 °C'
$ws.Range('E19').Value = '2026-02-15 17:49:06'
$ws.Range('O19').Value = '3.1 °C'
$ws.Range('E20').Value = '2026-02-15 17:49:07'
$ws.Range('O20').Value = '-3.0 °C'
$ws.Range('E21').Value = '2026-02-15 17:49:08'
$ws.Range('E22').Value = '2026-02-15 17:49:09'
$ws.Range('N22').Value = '-6.5 °C 17:18 TU'
$ws.Range('E23').Value = '2026-02-15 17:49:12'
$ws.Range('I23').Value = '1.5 mm'
$ws.Range('E24').Value = '2026-02-15 17:49:14'
$ws.Range('O24').Value = '8.6 °C'
$ws.Range('E25').Value = '2026-02-15 17:49:17'
$ws.Range('H25').Value = '61%'
$ws.Range('O25').Value = '-1.9 °C'
$ws.Range('E26').Value = '2026-02-15 17:49:19'
$ws.Range('E27').Value = '2026-02-15 17:49:21'
$ws.Range('H27').Value = '46%'
$ws.Range('E28').Value = '2026-02-15 17:49:24'
$ws.Range('H28').Value = '56%'
$ws.Range('O28').Value = '6.5 °C'
$ws.Range('E29').Value = '2026-02-15 17:49:27'
$ws.Range('H29').Value = '54%'
$ws.Range('E30').Value = '2026-02-15 17:49:29'
$ws.Range('H30').Value = '51%'
$ws.Range('O30').Value = '10.0 °C'
$ws.Range('E31').Value = '2026-02-15 17:49:32'
$ws.Range('O31').Value = '9.7 °C'
$ws.Range('E32').Value = '2026-02-15 17:49:34'
$ws.Range('H32').Value = '84%'
$ws.Range('K32').Value = '9.4 MJ/m2'
$ws.Range('O32').Value = '3.3 °C'
$ws.Range('E33').Value = '2026-02-15 17:49:37'
$ws.Range('H33').Value = '40%'
$ws.Range('O33').Value = '5.5 °C'
$ws.Range('E34').Value = '2026-02-15 17:49:39'
$ws.Range('O34').Value = '0.7 °C'
$ws.Range('E35').Value = '2026-02-15 17:49:42'
$ws.Range('O35').Value = '3.9 °C'
$ws.Range('E36').Value = '2026-02-15 17:49:44'
$ws.Range('H36').Value = '45%'
$ws.Range('E37').Value = '2026-02-15 17:49:47'
$ws.Range('J37').Value = '1016.2 hPa'
$ws.Range('E38').Value = '2026-02-15 17:49:49'
$ws.Range('O38').Value = '7.6 °C'
$ws.Range('E39').Value = '2026-02-15 17:49:52'
$ws.Range('H39').Value = '57%'
$ws.Range('O39').Value = '-3.3 °C'
$ws.Range('E40').Value = '2026-02-15 17:49:54'
$ws.Range('E41').Value = '2026-02-15 17:49:57'
$ws.Range('J41').Value = '1016.4 hPa'
$ws.Range('K41').Value = '12.5 MJ/m2'
$ws.Range('O41').Value = '12.2 °C'
$ws.Range('E42').Value = '2026-02-15 17:49:59'
$ws.Range('H42').Value = '53%'
$ws.Range('E43').Value = '2026-02-15 17:50:02'
$ws.Range('K43').Value = '12.9 MJ/m2'
$ws.Range('O43').Value = '6.0 °C'
$ws.Range('E44').Value = '2026-02-15 17:50:04'
$ws.Range('K44').Value = '9.8 MJ/m2'
$ws.Range('O44').Value = '-4.4 °C'
$ws.Range('E45').Value = '2026-02-15 17:50:07'
$ws.Range('H45').Value = '89%'
$ws.Range('J45').Value = '1023.7 hPa'
$ws.Range('O45').Value = '0.7 °C'
$ws.Range('E46').Value = '2026-02-15 17:50:09'
$ws.Range('J46').Value = '1019.6 hPa'
$ws.Range('K46').Value = '12.5 MJ/m2'
$ws.Range('O46').Value = '11.5 °C'
